$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4254.8965
$ws.Range("I132").Value = 3279.4167
$ws.Range("J132").Value = 8937.200000000001
$ws.Range("K132").Value = 9838.250100000001
$ws.Range("L132").Value = 26811.6
$ws.Range("M132").Value = -7308.250100000001
$ws.Range("N132").Value = -31871.6
$ws.Range("H138").Value = 2406.3027
$ws.Range("J138").Value = 2287.6726
$ws.Range("L138").Value = 6863.0178
$ws.Range("N138").Value = -17143.0178

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 32598
$ws.Range("J35").Value = 32598
$ws.Range("L35").Value = 32598
$ws.Range("N35").Value = -33218

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9180
$ws.Range("J50").Value = 9180
$ws.Range("L50").Value = 9180
$ws.Range("N50").Value = -10430
$ws.Range("H60").Value = 24584.428
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 26090.924
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 26090.924
$ws.Range("N60").Value = -27112.924
$ws.Range("M60").Value = -4489
$ws.Range("H68").Value = 17447.5
$ws.Range("J68").Value = 17447.5
$ws.Range("L68").Value = 17447.5
$ws.Range("N68").Value = -18945.5
$ws.Range("H71").Value = 17447.5
$ws.Range("J71").Value = 17447.5
$ws.Range("L71").Value = 52342.5
$ws.Range("N71").Value = -59830.5
$ws.Range("H99").Value = 1530.762
$ws.Range("I99").Value = 1535.6364
$ws.Range("J99").Value = 1525.4
$ws.Range("K99").Value = 1535.6364
$ws.Range("L99").Value = 1525.4
$ws.Range("M99").Value = -37.63640000000009
$ws.Range("N99").Value = -4521.4
$ws.Range("H126").Value = 1530.762
$ws.Range("I126").Value = 1535.6364
$ws.Range("J126").Value = 1525.4
$ws.Range("K126").Value = 4606.9092
$ws.Range("L126").Value = 4576.200000000001
$ws.Range("M126").Value = -2136.9092
$ws.Range("N126").Value = -9516.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 315.7143
$ws.Range("I97").Value = 185
$ws.Range("J97").Value = 368
$ws.Range("K97").Value = 555
$ws.Range("L97").Value = 1104
$ws.Range("M97").Value = -59
$ws.Range("N97").Value = -2096
$ws.Range("H131").Value = 931.2062
$ws.Range("I131").Value = 315
$ws.Range("J131").Value = 944.17896
$ws.Range("K131").Value = 945
$ws.Range("L131").Value = 2832.53688
$ws.Range("M131").Value = 4095
$ws.Range("N131").Value = -12912.53688

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 37696.816
$ws.Range("J123").Value = 37696.816
$ws.Range("L123").Value = 37696.816
$ws.Range("N123").Value = -42596.816
$ws.Range("H132").Value = 2530.36
$ws.Range("I132").Value = 2429.7222
$ws.Range("J132").Value = 2789.1428
$ws.Range("K132").Value = 7289.1666
$ws.Range("L132").Value = 8367.428400000001
$ws.Range("M132").Value = -4759.1666
$ws.Range("N132").Value = -13427.4284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5884584
$ws.Range("I7").Value = 7693616.5
$ws.Range("J7").Value = 5227.5
$ws.Range("K7").Value = 7693616.5
$ws.Range("L7").Value = 5227.5
$ws.Range("M7").Value = -7693504.5
$ws.Range("N7").Value = -5451.5
$ws.Range("H61").Value = 1964.3889
$ws.Range("I61").Value = 1304.1538
$ws.Range("J61").Value = 3681
$ws.Range("K61").Value = 1304.1538
$ws.Range("L61").Value = 3681
$ws.Range("M61").Value = -1102.1538
$ws.Range("N61").Value = -4085
$ws.Range("H100").Value = 1978.4445
$ws.Range("I100").Value = 1515.1428
$ws.Range("K100").Value = 1515.1428
$ws.Range("M100").Value = -974.1428000000001
$ws.Range("H104").Value = 17238.625
$ws.Range("J104").Value = 17238.625
$ws.Range("L104").Value = 17238.625
$ws.Range("N104").Value = -24226.625
$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680
$ws.Range("H109").Value = 24330
$ws.Range("J109").Value = 24330
$ws.Range("L109").Value = 24330
$ws.Range("N109").Value = -27104
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 54000
$ws.Range("J111").Value = 54000
$ws.Range("L111").Value = 54000
$ws.Range("N111").Value = -62180
$ws.Range("H113").Value = 1964.3889
$ws.Range("I113").Value = 1304.1538
$ws.Range("J113").Value = 3681
$ws.Range("K113").Value = 1304.1538
$ws.Range("L113").Value = 3681
$ws.Range("M113").Value = 865.8462
$ws.Range("N113").Value = -8021
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678
$ws.Range("H116").Value = 42250
$ws.Range("J116").Value = 42250
$ws.Range("L116").Value = 42250
$ws.Range("N116").Value = -51428
$ws.Range("H117").Value = 32000
$ws.Range("J117").Value = 32000
$ws.Range("L117").Value = 32000
$ws.Range("N117").Value = -41178
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 4933.3335
$ws.Range("I122").Value = 4166.6665
$ws.Range("J122").Value = 5316.6665
$ws.Range("K122").Value = 12499.9995
$ws.Range("L122").Value = 15949.9995
$ws.Range("M122").Value = -10049.9995
$ws.Range("N122").Value = -20849.9995
$ws.Range("H123").Value = 47492
$ws.Range("J123").Value = 47492
$ws.Range("L123").Value = 47492
$ws.Range("N123").Value = -57292
$ws.Range("H126").Value = 5884584
$ws.Range("I126").Value = 7693616.5
$ws.Range("J126").Value = 5227.5
$ws.Range("K126").Value = 23080849.5
$ws.Range("L126").Value = 15682.5
$ws.Range("M126").Value = -23078379.5
$ws.Range("N126").Value = -20622.5
$ws.Range("H127").Value = 54980
$ws.Range("J127").Value = 54980
$ws.Range("L127").Value = 54980
$ws.Range("N127").Value = -64900
$ws.Range("H128").Value = 47490
$ws.Range("J128").Value = 47490
$ws.Range("L128").Value = 47490
$ws.Range("N128").Value = -57450
$ws.Range("H129").Value = 40672
$ws.Range("J129").Value = 40672
$ws.Range("L129").Value = 40672
$ws.Range("N129").Value = -50672
$ws.Range("H130").Value = 51115
$ws.Range("J130").Value = 51115
$ws.Range("L130").Value = 51115
$ws.Range("N130").Value = -61155
$ws.Range("H131").Value = 37691.668
$ws.Range("J131").Value = 37691.668
$ws.Range("L131").Value = 37691.668
$ws.Range("N131").Value = -47771.668
$ws.Range("H132").Value = 5004.115
$ws.Range("I132").Value = 5101.0527
$ws.Range("J132").Value = 4741
$ws.Range("K132").Value = 15303.1581
$ws.Range("L132").Value = 14223
$ws.Range("M132").Value = -12773.1581
$ws.Range("N132").Value = -19283
$ws.Range("H136").Value = 1650.9487
$ws.Range("I136").Value = 1551
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4653
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2103
$ws.Range("N136").Value = -15600
